# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to match the values published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value  = 5067
$wsExpo.Range("F5").Value  = 7345
$wsExpo.Range("F12").Value = 4283
$wsExpo.Range("F13").Value = 1737
$wsExpo.Range("F16").Value = 2886
$wsExpo.Range("F21").Value = 426
$wsExpo.Range("F22").Value = 450
$wsExpo.Range("F23").Value = 290
$wsExpo.Range("F25").Value = 1679
$wsExpo.Range("F26").Value = 1163
$wsExpo.Range("F29").Value = 104
$wsExpo.Range("F32").Value = 510
$wsExpo.Range("F35").Value = 104
$wsExpo.Range("F37").Value = 2805
$wsExpo.Range("F39").Value = 25

# Sheet "全部类型" (sheet4) updates - same events, shifted one row down
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 5067
$wsAll.Range("F5").Value  = 7345
$wsAll.Range("F12").Value = 4283
$wsAll.Range("F13").Value = 1737
$wsAll.Range("F16").Value = 2886
$wsAll.Range("F21").Value = 426
$wsAll.Range("F22").Value = 450
$wsAll.Range("F23").Value = 290
$wsAll.Range("F25").Value = 1679
$wsAll.Range("F26").Value = 1163
$wsAll.Range("F29").Value = 104
$wsAll.Range("F32").Value = 510
$wsAll.Range("F35").Value = 104
$wsAll.Range("F37").Value = 2805
$wsAll.Range("F40").Value = 25
